$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Year values in column A for rows 2-5
$ws.Range("A2").Value = 2024
$ws.Range("A3").Value = 2021
$ws.Range("A4").Value = 2022
$ws.Range("A5").Value = 2400

# Update the active selection to A6 (mirrors the saved selection state)
$ws.Range("A6").Select()
